$p = $ppt.ActivePresentation
$hm = $p.HandoutMaster
$hm.HeadersFooters.Header.Text = "My handout header"
Write-Host "Done"
